$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 252. This shifts the existing rows 252-270
# down to 253-271, carrying their data (and the D-column date style) with
# them, matching the target diff where the old row 252 data reappears as
# row 253, old 253 as 254, ... old 270 as 271.
$ws.Rows("252:252").Insert()

# Populate the newly inserted row 252 with the new weekly record.
$ws.Range("A252").Value = 8
$ws.Range("B252").Value = "Terminal La Palmera de La Serena"
$ws.Range("C252").Value = "Coquimbo"
$ws.Range("D252").Value = 44610
$ws.Range("E252").Value = 4
$ws.Range("F252").Value = 100112032
$ws.Range("G252").Value = "Zapallo italiano"
$ws.Range("H252").Value = "Sin especificar"
$ws.Range("I252").Value = "Primera"
$ws.Range("J252").Value = 500
$ws.Range("K252").Value = 9000
$ws.Range("L252").Value = 10000
$ws.Range("M252").Value = 9500
$ws.Range("N252").Value = "`$/caja 60 unidades"
$ws.Range("O252").Value = "Provincia de Limarí"
$ws.Range("P252").Value = 158
$ws.Range("Q252").Value = 60
$ws.Range("R252").Value = "Hortaliza"

# Give the new row's date cell the same numeric date format used by the
# rest of column D (style index 2 in the original workbook).
$ws.Range("D252").NumberFormat = "YYYY-MM-DD HH:MM:SS"
